# Hjemme passive tweaks lichtwark deleted values
# - Row 1 (B1:E1): update the group-size header values.
# - Row 2 (B2:E2): B2 and D2 values are deleted outright; C2 and E2 get new values.
# - Row 3 (B3:E3): all four values are replaced with new values.
# - Reduce the selected range from the whole used range down to B1:E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 - B2 and D2 are removed entirely, C2 and E2 get new values
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -1.1077790549849604
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -0.7990084825247783

# Row 3 - all four values replaced
$ws.Range("B3").Value = -1.5242729628328515
$ws.Range("C3").Value = 0.2282749668256383
$ws.Range("D3").Value = -1.3531415117016905
$ws.Range("E3").Value = 1.9881402143597988

# Shrink the active selection to match the edited block
$ws.Range("B1:E3").Select()
